# The deck's "Integral" slide-master theme (ppt/theme/theme1.xml) is being
# swapped for the stock "Office Theme" palette (the theme previously only
# used by the notes master, ppt/theme/theme2.xml). Re-apply the 12 Office
# Theme scheme colors (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) onto
# the presentation's theme color scheme, which is what PowerPoint does under
# the hood when a different color theme is picked for the deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# RGB() values use the 0xBBGGRR packed format PowerPoint's COM layer expects.
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeThemeColors[$i - 1]
}
